# Weekly update: insert a new daily price record for
# "Pepino ensalada" at Vega Central Mapocho de Santiago, pushing the
# existing historical rows (444-467) down by one (to 445-468) and
# inserting the new week's data at row 444.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 444; this shifts rows 444:467 down to 445:468
# and extends the used range accordingly.
$ws.Rows(444).Insert()

# Populate the newly inserted row 444 with this week's record.
$ws.Cells.Item(444, 1).Value = 9
$ws.Cells.Item(444, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(444, 3).Value = "Metropolitana"
$ws.Cells.Item(444, 4).Value = 45267
$ws.Cells.Item(444, 5).Value = 13
$ws.Cells.Item(444, 6).Value = 100112043
$ws.Cells.Item(444, 7).Value = "Pepino ensalada"
$ws.Cells.Item(444, 8).Value = "Sin especificar"
$ws.Cells.Item(444, 9).Value = "Primera"
$ws.Cells.Item(444, 10).Value = 70
$ws.Cells.Item(444, 11).Value = 15000
$ws.Cells.Item(444, 12).Value = 16000
$ws.Cells.Item(444, 13).Value = 15500
$ws.Cells.Item(444, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(444, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(444, 16).Value = 258
$ws.Cells.Item(444, 17).Value = 60
$ws.Cells.Item(444, 18).Value = "Hortaliza"
